$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.435070037841797
$ws.Range("B1").Value = 1.895448565483093
$ws.Range("C1").Value = 2.130882263183594
$ws.Range("D1").Value = 2.43713903427124
$ws.Range("E1").Value = 2.968391418457031
